# The deck's applied Design/Theme ("Integral") is being switched back to the
# default "Office Theme" palette. The Integral theme's colours live in
# ppt/theme/theme2.xml -- the theme part actually applied to the slide
# master / presentation -- while ppt/theme/theme1.xml already holds the
# stock "Office Theme" palette (used by the notes master).
#
# PowerPoint's COM object model doesn't expose a "swap these two theme
# parts" operation directly, but re-colouring the applied theme's colour
# scheme to the Office Theme palette reproduces the same net effect on the
# presentation's live theme (ppt/theme/theme2.xml, the one the slide master
# / presentation actually renders with).

function HexToOleRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Stock "Office Theme" colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToOleRgb($officeThemeColors[$i - 1])
}
